# Add a new "canonical SMILES" column (D) that duplicates the existing
# "canonical isomeric SMILES" column (C) data, per the SAMPL6 microstates update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column (row 2, column D)
$ws.Cells.Item(2, 4).Value = "canonical SMILES"

# Copy each data row's canonical-isomeric-SMILES cell (column C) into the new
# column D, carrying over both the shared-string value and the row's cell
# formatting (fill/font/border/alignment), same as the existing C column.
for ($r = 3; $r -le 10; $r++) {
    $src = $ws.Cells.Item($r, 3)
    $dst = $ws.Cells.Item($r, 4)
    $src.Copy($dst)
}

# Match the authored column width for the new column D.
$ws.Columns.Item(4).ColumnWidth = 36
